$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A (username) updates
$ws.Range("A2").Value = "mngr322726"
$ws.Range("A3").Value = "admin"
$ws.Range("A4").Value = "mngr322726"
$ws.Range("A5").Value = "Eli"

# Column B (password) updates
$ws.Range("B2").Value = "ratahYn"
$ws.Range("B3").Value = "asfd"
$ws.Range("B4").Value = "ratahYn"
$ws.Range("B5").Value = "asdratahYn"

# Column C (Condition) updates
$ws.Range("C2").Value = "Y"
$ws.Range("C3").Value = "N"
$ws.Range("C4").Value = "N"
$ws.Range("C5").Value = "N"

# Update the active selection to B5 (matches saved view state)
$ws.Range("B5").Select()
